# Daily attendance processing
# Normalizes the "Recorded By" (column G) entries so the system account is
# listed after the human recorder's e-mail address, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$recordedByCol = 7  # Column G - "Recorded By"

$updated = 0
for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $recordedByCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Host "Reordered 'Recorded By' values in $updated row(s)."
